# Refresh cryptos list: updated coin order/links plus latest price & 1h volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.196.59'
$ws.Range("E2").Value = '  +1.64%  '

# Row 3
$ws.Range("D3").Value = '2.516.33'
$ws.Range("E3").Value = '  +1.30%  '

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").Value = "'590.61"
$ws.Range("E5").Value = '  +1.28%  '

# Row 6
$ws.Range("D6").Value = "'177.33"
$ws.Range("E6").Value = '  +3.64%  '

# Row 7
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("E8").Value = '  +0.99%  '

# Row 9
$ws.Range("D9").Value = "'0.146"
$ws.Range("E9").Value = '  +5.94%  '

# Row 10
$ws.Range("E10").Value = '  -0.49%  '

# Row 11
$ws.Range("D11").Value = "'0.341"
$ws.Range("E11").Value = '  +2.19%  '

# Row 12
$ws.Range("D12").Value = "'4.95"
$ws.Range("E12").Value = '  +0.63%  '

# Row 13
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.979.04'
$ws.Range("E13").Value = '  +2.30%  '

# Row 14
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = "'25.81"
$ws.Range("E14").Value = '  +1.75%  '

# Row 15
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '67.978.64'
$ws.Range("E15").Value = '  +1.48%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = "'0.0000173"
$ws.Range("E16").Value = '  +1.59%  '

# Row 17
$ws.Range("D17").Value = '2.501.68'
$ws.Range("E17").Value = '  -1.13%  '

# Row 18
$ws.Range("D18").Value = "'11.08"
$ws.Range("E18").Value = '  +0.85%  '

# Row 19
$ws.Range("D19").Value = "'7.58"
$ws.Range("E19").Value = '  +2.35%  '

# Row 20
$ws.Range("D20").Value = "'353.48"
$ws.Range("E20").Value = '  +1.44%  '

# Row 21
$ws.Range("E21").Value = '  +2.48%  '

# Row 22
$ws.Range("E22").Value = '  +0.10%  '

# Row 23
$ws.Range("D23").Value = "'71.08"
$ws.Range("E23").Value = '  +3.95%  '

# Row 24
$ws.Range("D24").Value = "'4.33"
$ws.Range("E24").Value = '  +2.54%  '

# Row 25
$ws.Range("D25").Value = "'1.77"
$ws.Range("E25").Value = '  -1.73%  '

# Row 26
$ws.Range("E26").Value = '  -0.96%  '

# Row 27
$ws.Range("D27").Value = '2.645.67'
$ws.Range("E27").Value = '  +1.27%  '

# Row 28
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = '  -0.08%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0924'
$ws.Range("E29").Value = '  +2.04%  '

# Row 30
$ws.Range("D30").Value = "'510.72"
$ws.Range("E30").Value = '  -0.33%  '

# Row 31
$ws.Range("D31").Value = "'7.84"
$ws.Range("E31").Value = '  +1.85%  '

# Row 32
$ws.Range("E32").Value = '  +3.25%  '

# Row 34
$ws.Range("E34").Value = '  +0.01%  '

# Row 35
$ws.Range("E35").Value = '  +3.60%  '

# Row 36
$ws.Range("D36").Value = "'165.01"
$ws.Range("E36").Value = '  +2.85%  '

# Row 37
$ws.Range("D37").Value = "'18.47"
$ws.Range("E37").Value = '  +1.48%  '

# Row 38
$ws.Range("D38").Value = "'18.66"
$ws.Range("E38").Value = '  -0.24%  '

# Row 39
$ws.Range("E39").Value = '  +0.80%  '

# Row 40
$ws.Range("E40").Value = '  +0.02%  '

# Row 41
$ws.Range("E41").Value = '  +3.54%  '

# Row 42
$ws.Range("E42").Value = '  +2.35%  '

# Row 43
$ws.Range("E43").Value = '  +0.84%  '

# Row 44
$ws.Range("E44").Value = '  +6.46%  '

# Row 45
$ws.Range("D45").Value = "'147.84"
$ws.Range("E45").Value = '  +3.50%  '

# Row 46
$ws.Range("D46").Value = "'3.55"
$ws.Range("E46").Value = '  +2.68%  '

# Row 47
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = "'0.522"
$ws.Range("E47").Value = '  +1.64%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0261'
$ws.Range("E48").Value = '  +3.96%  '

# Row 49
$ws.Range("D49").Value = "'0.0745"
$ws.Range("E49").Value = '  +2.25%  '

# Row 50
$ws.Range("E50").Value = '  +2.41%  '

# Row 51
$ws.Range("D51").Value = "'0.589"
$ws.Range("E51").Value = '  +1.30%  '
